$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9002.5
$ws.Range("J62").Value = 8670.333000000001
$ws.Range("L62").Value = 8670.333000000001
$ws.Range("N62").Value = -9918.333000000001
$ws.Range("H65").Value = 9002.5
$ws.Range("J65").Value = 8670.333000000001
$ws.Range("L65").Value = 43351.665
$ws.Range("N65").Value = -49591.665
$ws.Range("H112").Value = 3660.614
$ws.Range("J112").Value = 3721.0364
$ws.Range("L112").Value = 11163.1092
$ws.Range("N112").Value = -13379.1092
$ws.Range("H135").Value = 1764.12
$ws.Range("I135").Value = 1764.12
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 15877.08
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -13342.08
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1999.2894
$ws.Range("I137").Value = 1795.2069
$ws.Range("K137").Value = 5385.620699999999
$ws.Range("M137").Value = -2835.620699999999
$ws.Range("H138").Value = 3997.8965
$ws.Range("I138").Value = 2199.05
$ws.Range("J138").Value = 4944.6577
$ws.Range("K138").Value = 6597.150000000001
$ws.Range("L138").Value = 14833.9731
$ws.Range("M138").Value = -1457.150000000001
$ws.Range("N138").Value = -25113.9731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1154.5077
$ws.Range("I2").Value = 1035.619
$ws.Range("J2").Value = 4899.5
$ws.Range("K2").Value = 1035.619
$ws.Range("L2").Value = 4899.5
$ws.Range("M2").Value = -922.6189999999999
$ws.Range("N2").Value = -5125.5
$ws.Range("H61").Value = 4542.737
$ws.Range("I61").Value = 3600.8235
$ws.Range("K61").Value = 3600.8235
$ws.Range("M61").Value = -3388.8235
$ws.Range("H74").Value = 3644.111
$ws.Range("I74").Value = 3219.6
$ws.Range("J74").Value = 4174.75
$ws.Range("K74").Value = 3219.6
$ws.Range("L74").Value = 4174.75
$ws.Range("M74").Value = -2345.6
$ws.Range("N74").Value = -5922.75
$ws.Range("H77").Value = 3644.111
$ws.Range("I77").Value = 3219.6
$ws.Range("J77").Value = 4174.75
$ws.Range("K77").Value = 16098
$ws.Range("L77").Value = 20873.75
$ws.Range("M77").Value = -11730
$ws.Range("N77").Value = -29609.75
$ws.Range("H116").Value = 1154.5077
$ws.Range("I116").Value = 1035.619
$ws.Range("J116").Value = 4899.5
$ws.Range("K116").Value = 1035.619
$ws.Range("L116").Value = 4899.5
$ws.Range("M116").Value = 1258.381
$ws.Range("N116").Value = -9487.5
$ws.Range("H122").Value = 14128.765
$ws.Range("I122").Value = 16798.637
$ws.Range("J122").Value = 9234
$ws.Range("K122").Value = 50395.91099999999
$ws.Range("L122").Value = 27702
$ws.Range("M122").Value = -47945.91099999999
$ws.Range("N122").Value = -32602
$ws.Range("H136").Value = 4542.737
$ws.Range("I136").Value = 3600.8235
$ws.Range("K136").Value = 10802.4705
$ws.Range("M136").Value = -8252.470499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1154.5077
$ws.Range("I3").Value = 1035.619
$ws.Range("J3").Value = 4899.5
$ws.Range("K3").Value = 1035.619
$ws.Range("L3").Value = 4899.5
$ws.Range("M3").Value = -921.6189999999999
$ws.Range("N3").Value = -5127.5
$ws.Range("H80").Value = 1007.7273
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1007.7273
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1007.7273
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3003.7273
$ws.Range("H83").Value = 1007.7273
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1007.7273
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 5038.636500000001
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -15022.6365
$ws.Range("H86").Value = 16669117
$ws.Range("I86").Value = 20835584
$ws.Range("K86").Value = 20835584
$ws.Range("M86").Value = -20834461
$ws.Range("H89").Value = 16669117
$ws.Range("I89").Value = 20835584
$ws.Range("K89").Value = 104177920
$ws.Range("M89").Value = -104172304
$ws.Range("H105").Value = 4290.6
$ws.Range("I105").Value = 3865.25
$ws.Range("K105").Value = 3865.25
$ws.Range("M105").Value = -2118.25
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 30098.8
$ws.Range("I16").Value = 1497.5
$ws.Range("K16").Value = 1497.5
$ws.Range("M16").Value = -1210.5
$ws.Range("H31").Value = 7563.6387
$ws.Range("I31").Value = 7015.3125
$ws.Range("J31").Value = 8002.3
$ws.Range("K31").Value = 7015.3125
$ws.Range("L31").Value = 8002.3
$ws.Range("M31").Value = -6720.3125
$ws.Range("N31").Value = -8592.299999999999
$ws.Range("H34").Value = 7563.6387
$ws.Range("I34").Value = 7015.3125
$ws.Range("J34").Value = 8002.3
$ws.Range("K34").Value = 7015.3125
$ws.Range("L34").Value = 8002.3
$ws.Range("M34").Value = -6813.3125
$ws.Range("N34").Value = -8406.299999999999
$ws.Range("H113").Value = 30098.8
$ws.Range("I113").Value = 1497.5
$ws.Range("K113").Value = 1497.5
$ws.Range("M113").Value = 672.5
$ws.Range("H122").Value = 129948.52
$ws.Range("I122").Value = 166870.42
$ws.Range("K122").Value = 500611.26
$ws.Range("M122").Value = -498161.26
$ws.Range("H141").Value = 686566.7
$ws.Range("J141").Value = 686566.7
$ws.Range("L141").Value = 686566.7
$ws.Range("N141").Value = -696926.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 249750270
$ws.Range("I8").Value = 249750270
$ws.Range("K8").Value = 749250810
$ws.Range("M8").Value = -749250671
$ws.Range("H120").Value = 9103.223
$ws.Range("I120").Value = 9103.223
$ws.Range("K120").Value = 27309.669
$ws.Range("M120").Value = -22471.669

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 250
$ws.Range("I13").Value = 250
$ws.Range("K13").Value = 250
$ws.Range("M13").Value = -111
$ws.Range("H18").Value = 6006666
$ws.Range("I18").Value = 1509999.5
$ws.Range("J18").Value = 14999999
$ws.Range("K18").Value = 1509999.5
$ws.Range("L18").Value = 14999999
$ws.Range("M18").Value = -1509706.5
$ws.Range("N18").Value = -15000585
$ws.Range("H21").Value = 25000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 25000
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H35").Value = 17338.334
$ws.Range("I35").Value = 13507.5
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 13507.5
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -13209.5
$ws.Range("N35").Value = -25596
$ws.Range("H43").Value = 1823.4445
$ws.Range("I43").Value = 1823.4445
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1823.4445
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1672.4445
$ws.Range("N43").ClearContents()
$ws.Range("H70").Value = 5375
$ws.Range("J70").Value = 5375
$ws.Range("L70").Value = 5375
$ws.Range("N70").Value = -5915
$ws.Range("H73").Value = 5375
$ws.Range("J73").Value = 5375
$ws.Range("L73").Value = 5375
$ws.Range("N73").Value = -7247
$ws.Range("H132").Value = 8842.805
$ws.Range("I132").Value = 8121.2812
$ws.Range("K132").Value = 24363.8436
$ws.Range("M132").Value = -21833.8436

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4999.5
$ws.Range("I7").Value = 7000
$ws.Range("J7").Value = 4332.6665
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 4332.6665
$ws.Range("M7").Value = -6888
$ws.Range("N7").Value = -4556.6665
$ws.Range("H20").Value = 12500000
$ws.Range("I20").Value = 12500000
$ws.Range("K20").Value = 12500000
$ws.Range("M20").Value = -12499774
$ws.Range("H22").Value = 1355.3077
$ws.Range("I22").Value = 913.44446
$ws.Range("J22").Value = 2349.5
$ws.Range("K22").Value = 913.44446
$ws.Range("L22").Value = 2349.5
$ws.Range("M22").Value = -618.44446
$ws.Range("N22").Value = -2939.5
$ws.Range("H27").Value = 1355.3077
$ws.Range("I27").Value = 913.44446
$ws.Range("J27").Value = 2349.5
$ws.Range("K27").Value = 913.44446
$ws.Range("L27").Value = 2349.5
$ws.Range("M27").Value = -806.44446
$ws.Range("N27").Value = -2563.5
$ws.Range("H43").Value = 14946196
$ws.Range("I43").Value = 7107142
$ws.Range("K43").Value = 7107142
$ws.Range("M43").Value = -7106949
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H122").Value = 4713.2856
$ws.Range("J122").Value = 3502.5
$ws.Range("L122").Value = 10507.5
$ws.Range("N122").Value = -15407.5
$ws.Range("H126").Value = 4999.5
$ws.Range("I126").Value = 7000
$ws.Range("J126").Value = 4332.6665
$ws.Range("K126").Value = 21000
$ws.Range("L126").Value = 12997.9995
$ws.Range("M126").Value = -18530
$ws.Range("N126").Value = -17937.9995
$ws.Range("H132").Value = 2855.818
$ws.Range("I132").Value = 2166.2856
$ws.Range("K132").Value = 6498.8568
$ws.Range("M132").Value = -3968.8568
$ws.Range("H136").Value = 3573.9644
$ws.Range("I136").Value = 3368.1155
$ws.Range("K136").Value = 10104.3465
$ws.Range("M136").Value = -7554.3465

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4378.9287
$ws.Range("I122").Value = 2200.111
$ws.Range("K122").Value = 6600.333
$ws.Range("M122").Value = -4150.333
$ws.Range("H132").Value = 9142.054
$ws.Range("I132").Value = 7673.125
$ws.Range("K132").Value = 23019.375
$ws.Range("M132").Value = -20489.375
